$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: job 4082263309 - Frontend Developer @ PLAYA3ULL GAMES ---
$ws.Cells.Item(2,1).Value  = "'4082263309"
$ws.Cells.Item(2,2).Value  = "Frontend Developer"
$ws.Cells.Item(2,3).Value  = "https://www.linkedin.com/jobs/view/4082263309"
$ws.Cells.Item(2,4).Value  = "CHFZ0O2Zz5HJxmdHgCyYOw=="
$ws.Cells.Item(2,5).Value  = "'393453653"
$ws.Cells.Item(2,6).Value  = "PLAYA3ULL GAMES"
$ws.Cells.Item(2,7).Value  = "Australia (Remote)"
$ws.Cells.Item(2,9).Value  = "2024-11-24 11:43:32 +0000 UTC"
$ws.Cells.Item(2,12).Value = "yes"
$ws.Cells.Item(2,14).Value = "2024-11-25T08:12:03.032Z"
$ws.Cells.Item(2,15).Value = "2024-11-25T08:39:42.219Z"

# --- Row 3: job 4081862436 - Full Stack Developer - Freelance @ Twine ---
$ws.Cells.Item(3,1).Value  = "'4081862436"
$ws.Cells.Item(3,2).Value  = "Full Stack Developer - Freelance"
$ws.Cells.Item(3,3).Value  = "https://www.linkedin.com/jobs/view/4081862436"
$ws.Cells.Item(3,4).Value  = "CHFZ0O2Zz5HJxmdHgCyYOw=="
$ws.Cells.Item(3,5).Value  = ""
$ws.Cells.Item(3,6).Value  = "Twine"
$ws.Cells.Item(3,7).Value  = "Australia (Remote)"
$ws.Cells.Item(3,8).Value  = "Contract"
$ws.Cells.Item(3,9).Value  = "2024-11-21 12:30:39 +0000 UTC"
$ws.Cells.Item(3,12).Value = "yes"
$ws.Cells.Item(3,14).Value = "2024-11-25T08:12:03.074Z"
$ws.Cells.Item(3,15).Value = "2024-11-25T08:39:45.641Z"

# --- Row 4: job 4082423281 - React Developer @ Renaissance InfoSystems ---
$ws.Cells.Item(4,1).Value  = "'4082423281"
$ws.Cells.Item(4,2).Value  = "React Developer"
$ws.Cells.Item(4,3).Value  = "https://www.linkedin.com/jobs/view/4082423281"
$ws.Cells.Item(4,4).Value  = "CHFZ0O2Zz5HJxmdHgCyYOw=="
$ws.Cells.Item(4,5).Value  = "'856942361"
$ws.Cells.Item(4,6).Value  = "Renaissance InfoSystems"
$ws.Cells.Item(4,7).Value  = "Sydney, New South Wales, Australia (Hybrid)"
$ws.Cells.Item(4,8).Value  = "Contract"
$ws.Cells.Item(4,9).Value  = "2024-11-25 01:04:25 +0000 UTC"
$ws.Cells.Item(4,12).Value = "yes"
$ws.Cells.Item(4,14).Value = "2024-11-25T08:12:02.989Z"
$ws.Cells.Item(4,15).Value = "2024-11-25T08:39:38.684Z"

# --- Remove the now-obsolete rows 5-7 (job listings dropped off in this run) ---
$ws.Range("A5:O7").Delete()
